# Week01 report update:
#  - "Pending Task" list text re-worded (drop ".png" suffixes, normalize
#    casing/underscores) for every student still pending.
#  - Students PPP008 (row 9) and PPF007 (row 26) finished their last
#    remaining task, so their row flips from Pending -> Completed
#    (task list cleared, status cell + name cell get the "Completed"
#    green highlight style).
#  - Student PPP019 (row 19) still has one task left: "dulingo_update".
#  - "Generated" timestamp footer refreshed.
#  - Pending-Task column narrowed now that the text is shorter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTaskList = "Git_Task, Index_File_Updation, create_Html_file_on_Name, dulingo_update"

# Rows whose pending-task text is simply re-worded (still Pending).
$pendingRows = @(3, 6, 7, 8, 11, 14, 15, 16, 17, 18, 20, 23, 24, 25)
foreach ($r in $pendingRows) {
    $ws.Cells.Item($r, 4).Value = $newTaskList
}

# Row 19 keeps a single outstanding task.
$ws.Cells.Item(19, 4).Value = "dulingo_update"

# Rows 9 and 26: last pending task finished -> mark Completed.
$completeRows = @(9, 26)
foreach ($r in $completeRows) {
    $ws.Cells.Item($r, 4).ClearContents()
    $ws.Cells.Item($r, 5).Value = "Completed"

    # Match the workbook's existing "Completed" look (green fill / bold
    # white font) by copying the format from a row that already has it.
    $ws.Range("B2").Copy() | Out-Null
    $ws.Cells.Item($r, 2).PasteSpecial(-4122) | Out-Null
    $ws.Range("E2").Copy() | Out-Null
    $ws.Cells.Item($r, 5).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# Refresh the "Generated" timestamp footer.
$ws.Range("A29").Value = "Generated: 2023-09-02 11:47:58 AM"

# Pending-Task column is narrower now that the wording is shorter.
$ws.Columns.Item(4).ColumnWidth = 72.75
